$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 215; this shifts the existing rows
# 215..322 down to 216..323 (matching the dimension change A1:R322 -> A1:R323).
$ws.Rows.Item(215).Insert()

# Populate the newly inserted row 215 with its data.
$ws.Range("A215").Value = 9
$ws.Range("B215").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C215").Value = 'Metropolitana'
$ws.Range("D215").Value = 44806
$ws.Range("E215").Value = 13
$ws.Range("F215").Value = 100112017
$ws.Range("G215").Value = 'Apio'
$ws.Range("H215").Value = 'Americana (o)'
$ws.Range("I215").Value = 'Primera'
$ws.Range("J215").Value = 115
$ws.Range("K215").Value = 7000
$ws.Range("L215").Value = 10000
$ws.Range("M215").Value = 8652
$ws.Range("N215").Value = '$/docena de matas'
$ws.Range("O215").Value = 'Provincia del Elquí'
$ws.Range("P215").Value = 1442
$ws.Range("Q215").Value = 6
$ws.Range("R215").Value = 'Hortaliza'

# Make sure the D215 cell keeps the date/time number format used by the
# rest of the D column (style index 2 in the original file).
$ws.Range("D215").NumberFormat = $ws.Range("D216").NumberFormat
